$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new SceneQuest ("stonedoor2") to the QuestDungeon list for the
# --- "ancient tower" dungeon row (L6 in the Dungeon table). ---
$cell = $ws.Range("L6")
$current = $cell.Value()
$cell.Value = $current + "|stonedoor2;1"

# --- Tweak the workbook's Light 1 (Background 1) theme colour. ---
$themeColors = $wb.Theme.ThemeColorScheme
$light1 = $themeColors.Colors(2)
$light1.RGB = 13430215
